$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add "Customer Name" header in H1, copying the header style from G1 ---
$ws.Range("G1").Copy() | Out-Null
$ws.Range("H1").PasteSpecial(-4122) | Out-Null
$ws.Range("H1").Value = 'Customer Name'

# --- Fill column H ("Customer Name") for existing rows 2-26 ---
for ($r = 2; $r -le 26; $r++) {
    $ws.Cells.Item($r, 8).Value = "nan"
}

# --- Append new order rows 27-46 (columns A-H) ---
# Column A holds long numeric-looking Order IDs that must stay text (matches
# the existing rows inlineStr type) -- pre-format as Text so the COM layer
# does not silently coerce them to numbers.
$ws.Range("A27:A46").NumberFormat = "@"

# Row 27
$ws.Cells.Item(27, 1).Value = '20250308122815'
$ws.Cells.Item(27, 2).Value = '2025-03-08 12:28:15'
$ws.Cells.Item(27, 3).Value = '[{''Item Name'': ''Veg Thali'', ''Price'': 150.0}, {''Item Name'': ''Paneer Special'', ''Price'': 180.0}]'
$ws.Cells.Item(27, 4).Value = 330
$ws.Cells.Item(27, 5).Value = 'Not Paid'
$ws.Cells.Item(27, 6).Value = 'Take Away'
$ws.Cells.Item(27, 7).Value = 'Pending'
$ws.Cells.Item(27, 8).Value = 'nan'

# Row 28
$ws.Cells.Item(28, 1).Value = '20250308123025'
$ws.Cells.Item(28, 2).Value = '2025-03-08 12:30:25'
$ws.Cells.Item(28, 3).Value = '[{''Item Name'': ''Roti Sabzi'', ''Price'': 100.0}, {''Item Name'': ''Roti Sabzi'', ''Price'': 100.0}, {''Item Name'': ''Roti Sabzi'', ''Price'': 100.0}, {''Item Name'': ''Roti Sabzi'', ''Price'': 100.0}, {''Item Name'': ''South Indian Thali'', ''Price'': 200.0}, {''Item Name'': ''South Indian Thali'', ''Price'': 200.0}]'
$ws.Cells.Item(28, 4).Value = 800
$ws.Cells.Item(28, 5).Value = 'Paid'
$ws.Cells.Item(28, 6).Value = 'Take Away'
$ws.Cells.Item(28, 7).Value = 'Delivered'
$ws.Cells.Item(28, 8).Value = 'nan'

# Row 29
$ws.Cells.Item(29, 1).Value = '20250308124203'
$ws.Cells.Item(29, 2).Value = '2025-03-08 12:42:03'
$ws.Cells.Item(29, 3).Value = '[{''Item Name'': ''Veg Thali'', ''Price'': 150.0}, {''Item Name'': ''Veg Thali'', ''Price'': 150.0}, {''Item Name'': ''Veg Thali'', ''Price'': 150.0}]'
$ws.Cells.Item(29, 4).Value = 450
$ws.Cells.Item(29, 5).Value = 'Not Paid'
$ws.Cells.Item(29, 6).Value = 'Take Away'
$ws.Cells.Item(29, 7).Value = 'Pending'
$ws.Cells.Item(29, 8).Value = 'nan'

# Row 30
$ws.Cells.Item(30, 1).Value = '20250308131225'
$ws.Cells.Item(30, 2).Value = '2025-03-08 13:12:25'
$ws.Cells.Item(30, 3).Value = '[{''Item Name'': ''Roti Sabzi'', ''Price'': 100.0}]'
$ws.Cells.Item(30, 4).Value = 100
$ws.Cells.Item(30, 5).Value = 'Not Paid'
$ws.Cells.Item(30, 6).Value = 'Take Away'
$ws.Cells.Item(30, 7).Value = 'Pending'
$ws.Cells.Item(30, 8).Value = 'nan'

# Row 31
$ws.Cells.Item(31, 1).Value = '20250308131244'
$ws.Cells.Item(31, 2).Value = '2025-03-08 13:12:44'
$ws.Cells.Item(31, 3).Value = '[]'
$ws.Cells.Item(31, 4).Value = 0
$ws.Cells.Item(31, 5).Value = 'Not Paid'
$ws.Cells.Item(31, 6).Value = 'Take Away'
$ws.Cells.Item(31, 7).Value = 'Pending'
$ws.Cells.Item(31, 8).Value = 'nan'

# Row 32
$ws.Cells.Item(32, 1).Value = '20250308131304'
$ws.Cells.Item(32, 2).Value = '2025-03-08 13:13:04'
$ws.Cells.Item(32, 3).Value = '[{''Item Name'': ''Paneer Special'', ''Price'': 180.0}]'
$ws.Cells.Item(32, 4).Value = 180
$ws.Cells.Item(32, 5).Value = 'Not Paid'
$ws.Cells.Item(32, 6).Value = 'Take Away'
$ws.Cells.Item(32, 7).Value = 'Pending'
$ws.Cells.Item(32, 8).Value = 'nan'

# Row 33
$ws.Cells.Item(33, 1).Value = '20250308131424'
$ws.Cells.Item(33, 2).Value = '2025-03-08 13:14:24'
$ws.Cells.Item(33, 3).Value = '[{''Item Name'': ''Roti Sabzi'', ''Price'': 100.0}, {''Item Name'': ''Roti Sabzi'', ''Price'': 100.0}, {''Item Name'': ''Roti Sabzi'', ''Price'': 100.0}, {''Item Name'': ''Roti Sabzi'', ''Price'': 100.0}]'
$ws.Cells.Item(33, 4).Value = 400
$ws.Cells.Item(33, 5).Value = 'Not Paid'
$ws.Cells.Item(33, 6).Value = 'Take Away'
$ws.Cells.Item(33, 7).Value = 'Pending'
$ws.Cells.Item(33, 8).Value = 'aaaaaaaaaaaaa'

# Row 34
$ws.Cells.Item(34, 1).Value = '20250308131539'
$ws.Cells.Item(34, 2).Value = '2025-03-08 13:15:39'
$ws.Cells.Item(34, 3).Value = '[{''Item Name'': ''Roti Sabzi'', ''Price'': 100.0}, {''Item Name'': ''Roti Sabzi'', ''Price'': 100.0}, {''Item Name'': ''Roti Sabzi'', ''Price'': 100.0}]'
$ws.Cells.Item(34, 4).Value = 300
$ws.Cells.Item(34, 5).Value = 'Not Paid'
$ws.Cells.Item(34, 6).Value = 'Take Away'
$ws.Cells.Item(34, 7).Value = 'Pending'
$ws.Cells.Item(34, 8).Value = 'dfhkdjh'

# Row 35
$ws.Cells.Item(35, 1).Value = '20250308132414'
$ws.Cells.Item(35, 2).Value = '2025-03-08 13:24:14'
$ws.Cells.Item(35, 3).Value = '[{''Item Name'': ''Roti Sabzi'', ''Price'': 100.0}, {''Item Name'': ''South Indian Thali'', ''Price'': 200.0}]'
$ws.Cells.Item(35, 4).Value = 300
$ws.Cells.Item(35, 5).Value = 'Not Paid'
$ws.Cells.Item(35, 6).Value = 'Take Away'
$ws.Cells.Item(35, 7).Value = 'Pending'
$ws.Cells.Item(35, 8).Value = 'abc'

# Row 36
$ws.Cells.Item(36, 1).Value = '20250308132522'
$ws.Cells.Item(36, 2).Value = '2025-03-08 13:25:22'
$ws.Cells.Item(36, 3).Value = '[{''Item Name'': ''Veg Thali'', ''Price'': 150.0}, {''Item Name'': ''Veg Thali'', ''Price'': 150.0}, {''Item Name'': ''Veg Thali'', ''Price'': 150.0}]'
$ws.Cells.Item(36, 4).Value = 450
$ws.Cells.Item(36, 5).Value = 'Not Paid'
$ws.Cells.Item(36, 6).Value = 'Take Away'
$ws.Cells.Item(36, 7).Value = 'Pending'
$ws.Cells.Item(36, 8).Value = 'rajas'

# Row 37
$ws.Cells.Item(37, 1).Value = '20250308133046'
$ws.Cells.Item(37, 2).Value = '2025-03-08 13:30:46'
$ws.Cells.Item(37, 3).Value = '[{''Item Name'': ''Roti Sabzi'', ''Price'': 100.0}]'
$ws.Cells.Item(37, 4).Value = 100
$ws.Cells.Item(37, 5).Value = 'Not Paid'
$ws.Cells.Item(37, 6).Value = 'Take Away'
$ws.Cells.Item(37, 7).Value = 'Pending'
$ws.Cells.Item(37, 8).Value = 'asas'

# Row 38
$ws.Cells.Item(38, 1).Value = '20250308133209'
$ws.Cells.Item(38, 2).Value = '2025-03-08 13:32:09'
$ws.Cells.Item(38, 3).Value = '[{''Item Name'': ''Roti Sabzi'', ''Price'': 100.0}, {''Item Name'': ''Roti Sabzi'', ''Price'': 100.0}]'
$ws.Cells.Item(38, 4).Value = 200
$ws.Cells.Item(38, 5).Value = 'Not Paid'
$ws.Cells.Item(38, 6).Value = 'Take Away'
$ws.Cells.Item(38, 7).Value = 'Pending'
$ws.Cells.Item(38, 8).Value = 'nan'

# Row 39
$ws.Cells.Item(39, 1).Value = '20250308133215'
$ws.Cells.Item(39, 2).Value = '2025-03-08 13:32:15'
$ws.Cells.Item(39, 3).Value = '[{''Item Name'': ''Roti Sabzi'', ''Price'': 100.0}, {''Item Name'': ''Roti Sabzi'', ''Price'': 100.0}]'
$ws.Cells.Item(39, 4).Value = 200
$ws.Cells.Item(39, 5).Value = 'Not Paid'
$ws.Cells.Item(39, 6).Value = 'Take Away'
$ws.Cells.Item(39, 7).Value = 'Pending'
$ws.Cells.Item(39, 8).Value = 'nan'

# Row 40
$ws.Cells.Item(40, 1).Value = '20250308133215'
$ws.Cells.Item(40, 2).Value = '2025-03-08 13:32:15'
$ws.Cells.Item(40, 3).Value = '[{''Item Name'': ''Roti Sabzi'', ''Price'': 100.0}, {''Item Name'': ''Roti Sabzi'', ''Price'': 100.0}]'
$ws.Cells.Item(40, 4).Value = 200
$ws.Cells.Item(40, 5).Value = 'Not Paid'
$ws.Cells.Item(40, 6).Value = 'Take Away'
$ws.Cells.Item(40, 7).Value = 'Pending'
$ws.Cells.Item(40, 8).Value = 'nan'

# Row 41
$ws.Cells.Item(41, 1).Value = '20250308133215'
$ws.Cells.Item(41, 2).Value = '2025-03-08 13:32:15'
$ws.Cells.Item(41, 3).Value = '[{''Item Name'': ''Roti Sabzi'', ''Price'': 100.0}, {''Item Name'': ''Roti Sabzi'', ''Price'': 100.0}]'
$ws.Cells.Item(41, 4).Value = 200
$ws.Cells.Item(41, 5).Value = 'Not Paid'
$ws.Cells.Item(41, 6).Value = 'Take Away'
$ws.Cells.Item(41, 7).Value = 'Pending'
$ws.Cells.Item(41, 8).Value = 'nan'

# Row 42
$ws.Cells.Item(42, 1).Value = '20250308155055'
$ws.Cells.Item(42, 2).Value = '2025-03-08 15:50:55'
$ws.Cells.Item(42, 3).Value = '[{''Item Name'': ''Paneer Special'', ''Price'': 180.0}, {''Item Name'': ''Paneer Special'', ''Price'': 180.0}]'
$ws.Cells.Item(42, 4).Value = 360
$ws.Cells.Item(42, 5).Value = 'Not Paid'
$ws.Cells.Item(42, 6).Value = 'Take Away'
$ws.Cells.Item(42, 7).Value = 'Pending'
$ws.Cells.Item(42, 8).Value = 'nan'

# Row 43
$ws.Cells.Item(43, 1).Value = '20250308155314'
$ws.Cells.Item(43, 2).Value = '2025-03-08 15:53:14'
$ws.Cells.Item(43, 3).Value = '[{''Item Name'': ''Roti Sabzi'', ''Price'': 100.0}, {''Item Name'': ''Roti Sabzi'', ''Price'': 100.0}]'
$ws.Cells.Item(43, 4).Value = 200
$ws.Cells.Item(43, 5).Value = 'Not Paid'
$ws.Cells.Item(43, 6).Value = 'Take Away'
$ws.Cells.Item(43, 7).Value = 'Pending'
$ws.Cells.Item(43, 8).Value = 'nan'

# Row 44
$ws.Cells.Item(44, 1).Value = '20250308155322'
$ws.Cells.Item(44, 2).Value = '2025-03-08 15:53:22'
$ws.Cells.Item(44, 3).Value = '[{''Item Name'': ''Roti Sabzi'', ''Price'': 100.0}, {''Item Name'': ''Roti Sabzi'', ''Price'': 100.0}]'
$ws.Cells.Item(44, 4).Value = 200
$ws.Cells.Item(44, 5).Value = 'Not Paid'
$ws.Cells.Item(44, 6).Value = 'Take Away'
$ws.Cells.Item(44, 7).Value = 'Pending'
$ws.Cells.Item(44, 8).Value = 'sds'

# Row 45
$ws.Cells.Item(45, 1).Value = '20250308155356'
$ws.Cells.Item(45, 2).Value = '2025-03-08 15:53:56'
$ws.Cells.Item(45, 3).Value = '[{''Item Name'': ''South Indian Thali'', ''Price'': 200.0}, {''Item Name'': ''South Indian Thali'', ''Price'': 200.0}]'
$ws.Cells.Item(45, 4).Value = 400
$ws.Cells.Item(45, 5).Value = 'Not Paid'
$ws.Cells.Item(45, 6).Value = 'Take Away'
$ws.Cells.Item(45, 7).Value = 'Pending'
$ws.Cells.Item(45, 8).Value = 'nan'

# Row 46
$ws.Cells.Item(46, 1).Value = '20250308155419'
$ws.Cells.Item(46, 2).Value = '2025-03-08 15:54:19'
$ws.Cells.Item(46, 3).Value = '[{''Item Name'': ''Veg Thali'', ''Price'': 150.0}, {''Item Name'': ''Roti Sabzi'', ''Price'': 100.0}, {''Item Name'': ''Roti Sabzi'', ''Price'': 100.0}, {''Item Name'': ''Roti Sabzi'', ''Price'': 100.0}]'
$ws.Cells.Item(46, 4).Value = 450
$ws.Cells.Item(46, 5).Value = 'Not Paid'
$ws.Cells.Item(46, 6).Value = 'Take Away'
$ws.Cells.Item(46, 7).Value = 'Pending'
$ws.Cells.Item(46, 8).Value = 'rajas 1'

